$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# Update the panel query timestamps on the "data" sheet.
$data.Range("F2").Value = "2021-10-05 14:33:07.475639"
$data.Range("F3").Value = "2021-10-05 14:33:07.475647"

# Add the new "metadata" sheet right after "data".
$ws = $wb.Worksheets.Add($null, $data)
$ws.Name = "metadata"

# Reuse the bold/bordered header style from the "data" sheet's header row,
# and the "id column" style from its first data row.
$data.Range("B1:F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Header row.
$ws.Cells.Item(1, 2).Value = "data_name"
$ws.Cells.Item(1, 3).Value = "data_id"
$ws.Cells.Item(1, 4).Value = "data_version"
$ws.Cells.Item(1, 5).Value = "data_version_created"
$ws.Cells.Item(1, 6).Value = "panel_query_time"
$ws.Cells.Item(1, 7).Value = "panel_get_request"

# Data row.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Alagille syndrome"
$ws.Cells.Item(2, 3).Value = 36
$ws.Cells.Item(2, 4).Value = "'1.0"
$ws.Range("D2").Style = "Normal"
$ws.Cells.Item(2, 5).Value = "2020-09-24T11:24:48.443596Z"
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:33:07.471992"
$ws.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/36/?format=json"
